$d = $word.ActiveDocument

$d.Content.Find.Execute(
  "City’s 311 Center. See the codebook ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "City’s 311 Center. We have created a subset the data, limiting it to a sample of 1000 observations for the sake of faster loading for our demo purposes. See the codebook ",
  2)
